$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches the source file, which stores the
# "No." column as shared strings, not numbers) without leaving the cell's
# displayed style altered - restore the original style object afterwards.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

# 1. Rename the sheet.
$ws.Name = "BOM_Board1_LCD Add-on_2024-12-3"

# 2. Insert a new row at position 6, pushing the existing row 6
#    (MIC5504-3.3YM5-T5, "No."=5) and row 7 (trailing blank row) down to
#    rows 7 and 8 respectively.
$ws.Rows.Item(6).Insert()

# 3. Populate the newly inserted row 6 with the resistor BOM line.
Set-TextValue $ws.Cells.Item(6, 1) "5"
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(6, 3).Value = "10kΩ"
$ws.Cells.Item(6, 4).Value = "R1,R2"
$ws.Cells.Item(6, 5).Value = "R0805"
$ws.Cells.Item(6, 6).Value = "10kΩ"
$ws.Cells.Item(6, 7).Value = "0805W8F1002T5E"
$ws.Cells.Item(6, 8).Value = "UNI-ROYAL(厚声)"
$ws.Cells.Item(6, 9).Value = "C17414"
$ws.Cells.Item(6, 10).Value = "LCSC"

# 4. The previous row 6 (now row 7) keeps its contents, but its "No." bumps
#    from 5 to 6.
Set-TextValue $ws.Cells.Item(7, 1) "6"
